$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 490.18182
$ws.Range("I18").Value = 490.18182
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 490.18182
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -206.18182
$ws.Range("N18").ClearContents()
$ws.Range("H28").Value = 2555.7334
$ws.Range("I28").Value = 973.6
$ws.Range("K28").Value = 973.6
$ws.Range("M28").Value = -488.6
$ws.Range("H43").Value = 3373.8333
$ws.Range("J43").Value = 3747.6667
$ws.Range("L43").Value = 3747.6667
$ws.Range("N43").Value = -3885.6667
$ws.Range("H51").Value = 2403.923
$ws.Range("I51").Value = 1171
$ws.Range("J51").Value = 3174.5
$ws.Range("K51").Value = 1171
$ws.Range("L51").Value = 3174.5
$ws.Range("M51").Value = -687
$ws.Range("N51").Value = -4142.5
$ws.Range("H74").Value = 4704.476
$ws.Range("I74").Value = 3279.5
$ws.Range("K74").Value = 3279.5
$ws.Range("M74").Value = -2343.5
$ws.Range("H77").Value = 4704.476
$ws.Range("I77").Value = 3279.5
$ws.Range("K77").Value = 16397.5
$ws.Range("M77").Value = -11717.5
$ws.Range("H113").Value = 5019.5
$ws.Range("I113").Value = 3873
$ws.Range("J113").Value = 8459
$ws.Range("K113").Value = 3873
$ws.Range("L113").Value = 8459
$ws.Range("M113").Value = -619
$ws.Range("N113").Value = -14967
$ws.Range("H137").Value = 1973.1875
$ws.Range("I137").Value = 1894.1852
$ws.Range("K137").Value = 5682.5556
$ws.Range("M137").Value = -3132.5556
$ws.Range("H138").Value = 6671465.5
$ws.Range("I138").Value = 2087.3076
$ws.Range("J138").Value = 9014760
$ws.Range("K138").Value = 6261.9228
$ws.Range("L138").Value = 27044280
$ws.Range("M138").Value = -1121.9228
$ws.Range("N138").Value = -27054560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13168028
$ws.Range("I32").Value = 20837320
$ws.Range("J32").Value = 20668.535
$ws.Range("K32").Value = 20837320
$ws.Range("L32").Value = 20668.535
$ws.Range("M32").Value = -20837033
$ws.Range("N32").Value = -21242.535
$ws.Range("H33").Value = 14955
$ws.Range("I33").Value = 14955
$ws.Range("K33").Value = 14955
$ws.Range("M33").Value = -14626
$ws.Range("H61").Value = 20002762
$ws.Range("I61").Value = 24391876
$ws.Range("J61").Value = 7903.222
$ws.Range("K61").Value = 24391876
$ws.Range("L61").Value = 7903.222
$ws.Range("M61").Value = -24391664
$ws.Range("N61").Value = -8327.222
$ws.Range("H102").Value = 78283.08
$ws.Range("I102").Value = 81174.39999999999
$ws.Range("K102").Value = 81174.39999999999
$ws.Range("M102").Value = -79552.39999999999
$ws.Range("H110").Value = 15023.223
$ws.Range("J110").Value = 2977
$ws.Range("L110").Value = 2977
$ws.Range("N110").Value = -7067
$ws.Range("H132").Value = 24396666
$ws.Range("I132").Value = 7062.6665
$ws.Range("K132").Value = 21187.9995
$ws.Range("M132").Value = -18657.9995
$ws.Range("H136").Value = 20002762
$ws.Range("I136").Value = 24391876
$ws.Range("J136").Value = 7903.222
$ws.Range("K136").Value = 73175628
$ws.Range("L136").Value = 23709.666
$ws.Range("M136").Value = -73173078
$ws.Range("N136").Value = -28809.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2196.195
$ws.Range("I99").Value = 1367.1666
$ws.Range("K99").Value = 1367.1666
$ws.Range("M99").Value = 130.8334
$ws.Range("H105").Value = 11279.4
$ws.Range("I105").Value = 15272.571
$ws.Range("K105").Value = 15272.571
$ws.Range("M105").Value = -13525.571
$ws.Range("H107").Value = 3226.8845
$ws.Range("I107").Value = 2920.0588
$ws.Range("J107").Value = 3806.4443
$ws.Range("K107").Value = 2920.0588
$ws.Range("L107").Value = 3806.4443
$ws.Range("M107").Value = -1000.0588
$ws.Range("N107").Value = -7646.4443
$ws.Range("H132").Value = 114400
$ws.Range("J132").Value = 114400
$ws.Range("L132").Value = 114400
$ws.Range("N132").Value = -124520
$ws.Range("H134").Value = 2103.9443
$ws.Range("I134").Value = 2102.5103
$ws.Range("J134").Value = 2118
$ws.Range("K134").Value = 6307.5309
$ws.Range("L134").Value = 6354
$ws.Range("M134").Value = -3772.5309
$ws.Range("N134").Value = -11424

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1420.9412
$ws.Range("I134").Value = 1259.8125
$ws.Range("K134").Value = 3779.4375
$ws.Range("M134").Value = -1244.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 915.5333000000001
$ws.Range("I25").Value = 184
$ws.Range("J25").Value = 2378.6
$ws.Range("K25").Value = 552
$ws.Range("L25").Value = 7135.799999999999
$ws.Range("M25").Value = -383
$ws.Range("N25").Value = -7473.799999999999
$ws.Range("H30").Value = 915.5333000000001
$ws.Range("I30").Value = 184
$ws.Range("J30").Value = 2378.6
$ws.Range("K30").Value = 552
$ws.Range("L30").Value = 7135.799999999999
$ws.Range("M30").Value = -450
$ws.Range("N30").Value = -7339.799999999999
$ws.Range("H140").Value = 2211.4546
$ws.Range("I140").Value = 1589.2
$ws.Range("J140").Value = 2730
$ws.Range("K140").Value = 4767.6
$ws.Range("L140").Value = 8190
$ws.Range("M140").Value = 412.3999999999996
$ws.Range("N140").Value = -18550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2949.65
$ws.Range("I80").Value = 2826.3
$ws.Range("J80").Value = 3073
$ws.Range("K80").Value = 2826.3
$ws.Range("L80").Value = 3073
$ws.Range("M80").Value = -1828.3
$ws.Range("N80").Value = -5069
$ws.Range("H83").Value = 2949.65
$ws.Range("I83").Value = 2826.3
$ws.Range("J83").Value = 3073
$ws.Range("K83").Value = 14131.5
$ws.Range("L83").Value = 15365
$ws.Range("M83").Value = -9139.5
$ws.Range("N83").Value = -25349
$ws.Range("H98").Value = 386707.5
$ws.Range("J98").Value = 386707.5
$ws.Range("L98").Value = 386707.5
$ws.Range("N98").Value = -392697.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 490.8125
$ws.Range("I16").Value = 490.8125
$ws.Range("K16").Value = 490.8125
$ws.Range("M16").Value = -320.8125
$ws.Range("H22").Value = 2555.6287
$ws.Range("J22").Value = 3064.85
$ws.Range("L22").Value = 3064.85
$ws.Range("N22").Value = -3654.85
$ws.Range("H27").Value = 2555.6287
$ws.Range("J27").Value = 3064.85
$ws.Range("L27").Value = 3064.85
$ws.Range("N27").Value = -3278.85
$ws.Range("H40").Value = 4751.8887
$ws.Range("I40").Value = 4050.1
$ws.Range("J40").Value = 6757
$ws.Range("K40").Value = 4050.1
$ws.Range("L40").Value = 6757
$ws.Range("M40").Value = -3914.1
$ws.Range("N40").Value = -7029
$ws.Range("H61").Value = 1963.9412
$ws.Range("I61").Value = 781.3077
$ws.Range("K61").Value = 781.3077
$ws.Range("M61").Value = -579.3077
$ws.Range("H68").Value = 2591.0435
$ws.Range("J68").Value = 3547.5
$ws.Range("L68").Value = 3547.5
$ws.Range("N68").Value = -5045.5
$ws.Range("H71").Value = 2591.0435
$ws.Range("J71").Value = 3547.5
$ws.Range("L71").Value = 17737.5
$ws.Range("N71").Value = -25225.5
$ws.Range("H100").Value = 2095.4546
$ws.Range("I100").Value = 1456.125
$ws.Range("K100").Value = 1456.125
$ws.Range("M100").Value = -915.125
$ws.Range("H113").Value = 1963.9412
$ws.Range("I113").Value = 781.3077
$ws.Range("K113").Value = 781.3077
$ws.Range("M113").Value = 1388.6923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 968.75
$ws.Range("I81").Value = 625
$ws.Range("K81").Value = 1250
$ws.Range("M81").Value = -189
$ws.Range("H84").Value = 968.75
$ws.Range("I84").Value = 625
$ws.Range("K84").Value = 6250
$ws.Range("M84").Value = -946
$ws.Range("H92").Value = 38349.668
$ws.Range("J92").Value = 38349.668
$ws.Range("L92").Value = 38349.668
$ws.Range("N92").Value = -43341.668
$ws.Range("H113").Value = 394.85715
$ws.Range("I113").Value = 207.45454
$ws.Range("K113").Value = 622.3636200000001
$ws.Range("M113").Value = 1547.63638
$ws.Range("H132").Value = 4656.592
$ws.Range("I132").Value = 4664.7046
$ws.Range("K132").Value = 13994.1138
$ws.Range("M132").Value = -11464.1138
